$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "29.520.55"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "1.912.47"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("E4").Value = "  +0.74%  "
$ws.Range("D5").Value = "325.99"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("D7").Value = "0.4826"
$ws.Range("E7").Value = "  +1.63%  "
$ws.Range("D8").Value = "0.4070"
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("D9").Value = "0.08166"
$ws.Range("E9").Value = "  +1.38%  "
$ws.Range("D10").Value = "1.012"
$ws.Range("E10").Value = "  +0.22%  "
$ws.Range("D11").Value = "23.36"
$ws.Range("E11").Value = "  +4.15%  "
$ws.Range("D12").Value = "1.918.64"
$ws.Range("E12").Value = "  +3.37%  "
$ws.Range("D13").Value = "6.000"
$ws.Range("E13").Value = "  +1.24%  "
$ws.Range("D14").Value = "7.121"
$ws.Range("E14").Value = "  -0.41%  "
$ws.Range("D15").Value = "90.37"
$ws.Range("E15").Value = "  +0.95%  "
$ws.Range("D16").Value = "0.06795"
$ws.Range("E16").Value = "  +2.88%  "
$ws.Range("E17").Value = "  +0.80%  "
$ws.Range("D18").Value = "0.00001041"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("D19").Value = "17.70"
$ws.Range("E19").Value = "  -0.16%  "
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("D21").Value = "29.530.75"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").Value = "5.616"
$ws.Range("E22").Value = "  +1.40%  "
$ws.Range("D23").Value = "11.81"
$ws.Range("E23").Value = "  +2.64%  "
$ws.Range("E24").Value = "  -1.15%  "
$ws.Range("D25").Value = "2.144.62"
$ws.Range("E25").Value = "  +2.14%  "
$ws.Range("D26").Value = "155.63"
$ws.Range("E26").Value = "  +0.66%  "
$ws.Range("D27").Value = "6.413"
$ws.Range("E27").Value = "  +8.82%  "
$ws.Range("D28").Value = "20.09"
$ws.Range("E28").Value = "  +1.40%  "
$ws.Range("D29").Value = "2.097"
$ws.Range("E29").Value = "  -1.62%  "
$ws.Range("D30").Value = "119.66"
$ws.Range("E30").Value = "  +1.73%  "
$ws.Range("D31").Value = "1.029"
$ws.Range("E31").Value = "  -2.23%  "
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("D33").Value = "5.525"
$ws.Range("E33").Value = "  +2.26%  "
$ws.Range("D34").Value = "3.559"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("E35").Value = "  -2.38%  "
$ws.Range("E36").Value = "  +0.33%  "
$ws.Range("D37").Value = "0.06091"
$ws.Range("E37").Value = "  -0.45%  "
$ws.Range("D38").Value = "1.182"
$ws.Range("E38").Value = "  +0.73%  "
$ws.Range("D39").Value = "10.79"
$ws.Range("E39").Value = "  +6.31%  "
$ws.Range("D40").Value = "0.5939"
$ws.Range("E40").Value = "  +0.93%  "
$ws.Range("D41").Value = "7.944"
$ws.Range("E41").Value = "  -4.44%  "
$ws.Range("D42").Value = "0.1859"
$ws.Range("E42").Value = "  +0.80%  "
$ws.Range("D43").Value = "2.481"
$ws.Range("E43").Value = "  -2.40%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "0.07727"
$ws.Range("E45").Value = "  -3.58%  "
$ws.Range("D46").Value = "12.45"
$ws.Range("E46").Value = "  +2.12%  "
$ws.Range("D47").Value = "0.5568"
$ws.Range("E47").Value = "  +0.43%  "
$ws.Range("D48").Value = "1.944"
$ws.Range("E48").Value = "  +0.43%  "
$ws.Range("D49").Value = "115.80"
$ws.Range("E49").Value = "  +2.35%  "
$ws.Range("D50").Value = "72.66"
$ws.Range("E50").Value = "  +1.53%  "
$ws.Range("E51").Value = "  +1.99%  "

$ws.Range("D5").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D18").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D28").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D31").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D50").ClearFormats()
